$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new progress-report answers (adds 4 new shared strings).
$ws.Range("E2").Value = "ik ben nu aan het kijken om het ""hoofd"", lijnsensor en addon kaart te vernieuwen maar ik ga waschijnlijk niet het hoofdbord kunnen vervangen maar meschien al tekenen"
$ws.Range("F4").Value = "door te expirimenteren met wifi kan ik het op afstand besturen maar school maakt het wel moeilijker"
$ws.Range("C7").Value = "de nieuwe bordjes bestellen maken en testen en mocht ik nog tijd hebben al beginnen met het hoofdbord te tekenen"
$ws.Range("F5").Value = "dit vind ik heel moeilijk maar ik werk wel veet thuis"

# Row heights shrink now that the new text has been entered (auto sizing).
$ws.Rows.Item(2).RowHeight = 72.5
$ws.Rows.Item(3).RowHeight = 43.5
$ws.Rows.Item(4).RowHeight = 43.5

# Leave the selection where the user was last typing.
$ws.Range("F5").Select() | Out-Null
